$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PB Items")
$ws.Activate()

# Fix a typo in the Priority column for US4 ("Must have" -> "Must Have")
$ws.Range("F9").Value = "Must Have"

# Fill in the Status column (H) for every backlog item: "To be started"
# Copy formatting from the already-styled neighbour (I6) so the new cells
# pick up the same cell style (s="1") used throughout row 6 / column I.
$ws.Range("I6").Copy() | Out-Null
$statusRange = $ws.Range("H6:H14")
$statusRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("H6").Value = "To be started"
$ws.Range("H7").Value = "To be started"
$ws.Range("H8").Value = "To be started"
$ws.Range("H9").Value = "To be started"
$ws.Range("H10").Value = "To be started"
$ws.Range("H11").Value = "To be started"
$ws.Range("H12").Value = "To be started"
$ws.Range("H13").Value = "To be started"
$ws.Range("H14").Value = "To be started"

# Fill in the Estimation column (J) with Fibonacci-series story points
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 6
$ws.Range("J8").Value = 4
$ws.Range("J9").Value = 3
$ws.Range("J10").Value = 2
$ws.Range("J11").Value = 3
$ws.Range("J12").Value = 2
$ws.Range("J13").Value = 4
$ws.Range("J14").Value = 5

# Resize columns: split the old merged G:H-width column 7-8 so column 8
# (Status) is wider, and narrow column 10 (Estimation) since it now only
# holds single-digit numbers.
$ws.Columns.Item(8).ColumnWidth = 14.8
$ws.Columns.Item(10).ColumnWidth = 6.14

# Update the active selection / scroll position to match the saved view
$ws.Range("G17").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
